$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "MuSCs"
$ws.Cells.Item(2, 2).Value = "Calca"
$ws.Cells.Item(2, 3).Value = "Ramp1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.5
$ws.Cells.Item(2, 7).Value = 0.354751
$ws.Cells.Item(2, 8).Value = 0.709502
$ws.Cells.Item(2, 9).Value = 0.7067029850439027
$ws.Cells.Item(2, 10).Value = 0.6163205031315422
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.22082
$ws.Cells.Item(2, 14).Value = 0.66246
$ws.Cells.Item(2, 15).Value = 0.03343653775821253
$ws.Cells.Item(2, 16).Value = 0.03487484177052758
$ws.Cells.Item(2, 17).Value = 0.07833611582
$ws.Cells.Item(2, 18).Value = 0.47001669492
$ws.Cells.Item(2, 19).Value = 0.02362970104326196
$ws.Cells.Item(2, 20).Value = 0.02149408002664448

# Row 3
$ws.Cells.Item(3, 1).Value = "MuSCs"
$ws.Cells.Item(3, 2).Value = "Calca"
$ws.Cells.Item(3, 3).Value = "Ramp1"
$ws.Cells.Item(3, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.5
$ws.Cells.Item(3, 7).Value = 0.354751
$ws.Cells.Item(3, 8).Value = 0.709502
$ws.Cells.Item(3, 9).Value = 0.7067029850439027
$ws.Cells.Item(3, 10).Value = 0.6163205031315422
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.616755666666667
$ws.Cells.Item(3, 14).Value = 4.850267000000001
$ws.Cells.Item(3, 15).Value = 0.244808947986161
$ws.Cells.Item(3, 16).Value = 0.2553396343474497
$ws.Cells.Item(3, 17).Value = 0.5735456895056668
$ws.Cells.Item(3, 18).Value = 3.441274137034
$ws.Cells.Item(3, 19).Value = 0.1730072143072775
$ws.Cells.Item(3, 20).Value = 0.1573710519104442

# Row 4
$ws.Cells.Item(4, 1).Value = "MuSCs"
$ws.Cells.Item(4, 2).Value = "Calca"
$ws.Cells.Item(4, 3).Value = "Ramp1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.5
$ws.Cells.Item(4, 7).Value = 0.354751
$ws.Cells.Item(4, 8).Value = 0.709502
$ws.Cells.Item(4, 9).Value = 0.7067029850439027
$ws.Cells.Item(4, 10).Value = 0.6163205031315422
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.817103
$ws.Cells.Item(4, 14).Value = 1.634206
$ws.Cells.Item(4, 15).Value = 0.1237256376770616
$ws.Cells.Item(4, 16).Value = 0.08603187463461461
$ws.Cells.Item(4, 17).Value = 0.289868106353
$ws.Cells.Item(4, 18).Value = 1.159472425412
$ws.Cells.Item(4, 19).Value = 0.08743727747283976
$ws.Cells.Item(4, 20).Value = 0.05302320826015543

# Row 5
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Calca"
$ws.Cells.Item(5, 3).Value = "Ramp1"
$ws.Cells.Item(5, 4).Value = "Neutrophils"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.5
$ws.Cells.Item(5, 7).Value = 0.354751
$ws.Cells.Item(5, 8).Value = 0.709502
$ws.Cells.Item(5, 9).Value = 0.7067029850439027
$ws.Cells.Item(5, 10).Value = 0.6163205031315422
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.532357333333334
$ws.Cells.Item(5, 14).Value = 7.597072000000001
$ws.Cells.Item(5, 15).Value = 0.38344924188609
$ws.Cells.Item(5, 16).Value = 0.3999436704394311
$ws.Cells.Item(5, 17).Value = 0.8983562963573334
$ws.Cells.Item(5, 18).Value = 5.390137778144
$ws.Cells.Item(5, 19).Value = 0.2709847238537213
$ws.Cells.Item(5, 20).Value = 0.2464934841895059

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Calca"
$ws.Cells.Item(6, 3).Value = "Ramp1"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.5
$ws.Cells.Item(6, 7).Value = 0.354751
$ws.Cells.Item(6, 8).Value = 0.709502
$ws.Cells.Item(6, 9).Value = 0.7067029850439027
$ws.Cells.Item(6, 10).Value = 0.6163205031315422
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.417116666666667
$ws.Cells.Item(6, 14).Value = 4.25135
$ws.Cells.Item(6, 15).Value = 0.2145796346924748
$ws.Cells.Item(6, 16).Value = 0.223809978807977
$ws.Cells.Item(6, 17).Value = 0.5027235546166667
$ws.Cells.Item(6, 18).Value = 3.0163413277
$ws.Cells.Item(6, 19).Value = 0.1516440683668021
$ws.Cells.Item(6, 20).Value = 0.1379386787447922

# Row 7
$ws.Cells.Item(7, 1).Value = "Neutrophils"
$ws.Cells.Item(7, 2).Value = "Calca"
$ws.Cells.Item(7, 3).Value = "Ramp1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.1472293333333334
$ws.Cells.Item(7, 8).Value = 0.441688
$ws.Cells.Item(7, 9).Value = 0.2932970149560972
$ws.Cells.Item(7, 10).Value = 0.3836794968684579
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.22082
$ws.Cells.Item(7, 14).Value = 0.66246
$ws.Cells.Item(7, 15).Value = 0.03343653775821253
$ws.Cells.Item(7, 16).Value = 0.03487484177052758
$ws.Cells.Item(7, 17).Value = 0.03251118138666667
$ws.Cells.Item(7, 18).Value = 0.29260063248
$ws.Cells.Item(7, 19).Value = 0.00980683671495057
$ws.Cells.Item(7, 20).Value = 0.0133807617438831

# Row 8
$ws.Cells.Item(8, 1).Value = "Neutrophils"
$ws.Cells.Item(8, 2).Value = "Calca"
$ws.Cells.Item(8, 3).Value = "Ramp1"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.1472293333333334
$ws.Cells.Item(8, 8).Value = 0.441688
$ws.Cells.Item(8, 9).Value = 0.2932970149560972
$ws.Cells.Item(8, 10).Value = 0.3836794968684579
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.616755666666667
$ws.Cells.Item(8, 14).Value = 4.850267000000001
$ws.Cells.Item(8, 15).Value = 0.244808947986161
$ws.Cells.Item(8, 16).Value = 0.2553396343474497
$ws.Cells.Item(8, 17).Value = 0.2380338589662223
$ws.Cells.Item(8, 18).Value = 2.142304730696
$ws.Cells.Item(8, 19).Value = 0.0718017336788835
$ws.Cells.Item(8, 20).Value = 0.0979685824370055

# Row 9
$ws.Cells.Item(9, 1).Value = "Neutrophils"
$ws.Cells.Item(9, 2).Value = "Calca"
$ws.Cells.Item(9, 3).Value = "Ramp1"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.1472293333333334
$ws.Cells.Item(9, 8).Value = 0.441688
$ws.Cells.Item(9, 9).Value = 0.2932970149560972
$ws.Cells.Item(9, 10).Value = 0.3836794968684579
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.817103
$ws.Cells.Item(9, 14).Value = 1.634206
$ws.Cells.Item(9, 15).Value = 0.1237256376770616
$ws.Cells.Item(9, 16).Value = 0.08603187463461461
$ws.Cells.Item(9, 17).Value = 0.1203015299546667
$ws.Cells.Item(9, 18).Value = 0.7218091797280001
$ws.Cells.Item(9, 19).Value = 0.03628836020422179
$ws.Cells.Item(9, 20).Value = 0.03300866637445918

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutrophils"
$ws.Cells.Item(10, 2).Value = "Calca"
$ws.Cells.Item(10, 3).Value = "Ramp1"
$ws.Cells.Item(10, 4).Value = "Neutrophils"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1472293333333334
$ws.Cells.Item(10, 8).Value = 0.441688
$ws.Cells.Item(10, 9).Value = 0.2932970149560972
$ws.Cells.Item(10, 10).Value = 0.3836794968684579
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.532357333333334
$ws.Cells.Item(10, 14).Value = 7.597072000000001
$ws.Cells.Item(10, 15).Value = 0.38344924188609
$ws.Cells.Item(10, 16).Value = 0.3999436704394311
$ws.Cells.Item(10, 17).Value = 0.3728372819484445
$ws.Cells.Item(10, 18).Value = 3.355535537536
$ws.Cells.Item(10, 19).Value = 0.1124645180323687
$ws.Cells.Item(10, 20).Value = 0.1534501862499253

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutrophils"
$ws.Cells.Item(11, 2).Value = "Calca"
$ws.Cells.Item(11, 3).Value = "Ramp1"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.1472293333333334
$ws.Cells.Item(11, 8).Value = 0.441688
$ws.Cells.Item(11, 9).Value = 0.2932970149560972
$ws.Cells.Item(11, 10).Value = 0.3836794968684579
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.417116666666667
$ws.Cells.Item(11, 14).Value = 4.25135
$ws.Cells.Item(11, 15).Value = 0.2145796346924748
$ws.Cells.Item(11, 16).Value = 0.223809978807977
$ws.Cells.Item(11, 17).Value = 0.2086411420888889
$ws.Cells.Item(11, 18).Value = 1.8777702788
$ws.Cells.Item(11, 19).Value = 0.06293556632567265
$ws.Cells.Item(11, 20).Value = 0.08587130006318484
